$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("main")

# Clear stray leftover values in rows 3-5 (unique/top/freq) for columns G:W,Y,AC,AG
$ws.Range("G3:W5").ClearContents()
$ws.Range("Y3:Y5").ClearContents()
$ws.Range("AC3:AC5").ClearContents()
$ws.Range("AG3:AG5").ClearContents()

# Populate descriptive-statistics values for rows 6-12 (mean,std,min,25%,50%,75%,max)

# Row 6
$ws.Range("G6").Value = 47.04947321487399
$ws.Range("H6").Value = 36.65452152153802
$ws.Range("I6").Value = 51.29760338579531
$ws.Range("J6").Value = 50.86137493920096
$ws.Range("K6").Value = 2237.825678954217
$ws.Range("L6").Value = 2015.78642389924
$ws.Range("M6").Value = 2781.261072589395
$ws.Range("N6").Value = 151185.0129286476
$ws.Range("O6").Value = 2.882929392896679
$ws.Range("P6").Value = 221.4736711023251
$ws.Range("Q6").Value = 28.67589288437186
$ws.Range("R6").Value = 1.221019912769642
$ws.Range("S6").Value = 0.7893919278945938
$ws.Range("T6").Value = 5.068735332600547
$ws.Range("U6").Value = 37.32600118667948
$ws.Range("V6").Value = 1158618760.272202
$ws.Range("W6").Value = 37.70113294315134
$ws.Range("Y6").Value = 0.1791045396782302
$ws.Range("AC6").Value = 35.30624285197852
$ws.Range("AG6").Value = 306.7955671989235

# Row 7
$ws.Range("G7").Value = 18.1601539559753
$ws.Range("H7").Value = 23.84502287773762
$ws.Range("I7").Value = 21.02170683863201
$ws.Range("J7").Value = 21.20635832269786
$ws.Range("K7").Value = 1938.027249409436
$ws.Range("L7").Value = 1751.787828013133
$ws.Range("M7").Value = 1875.330867057288
$ws.Range("N7").Value = 148477.7993157351
$ws.Range("O7").Value = 8.734527958605121
$ws.Range("P7").Value = 6927.875627478162
$ws.Range("Q7").Value = 72.35422130088232
$ws.Range("R7").Value = 0.7286885655475006
$ws.Range("S7").Value = 0.5456071293440071
$ws.Range("T7").Value = 9.323698036514696
$ws.Range("U7").Value = 21.75485270480808
$ws.Range("V7").Value = 2882812342.31101
$ws.Range("W7").Value = 24.89441352224999
$ws.Range("Y7").Value = 3.579519170148852
$ws.Range("AC7").Value = 1337.121788481756
$ws.Range("AG7").Value = 9900.813149134461

# Row 8
$ws.Range("G8").Value = 5.45525079259131
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 5.38395904436859
$ws.Range("J8").Value = 1.80598555211558
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = 29.11385274232588
$ws.Range("N8").Value = 162.3469609285522
$ws.Range("O8").Value = 0.01
$ws.Range("P8").Value = 0.196309744053503
$ws.Range("Q8").Value = 0.007169704000657
$ws.Range("R8").Value = -2.37335245719635
$ws.Range("S8").Value = -0.231835716417949
$ws.Range("T8").Value = -46.7852118737
$ws.Range("U8").Value = 9.590779098391289
$ws.Range("V8").Value = 4423728.81355932
$ws.Range("W8").Value = 8.961593506751379
$ws.Range("Y8").Value = 0.0006846846846849068
$ws.Range("AC8").Value = 0.048628716485443
$ws.Range("AG8").Value = 0.1107983699353448

# Row 9
$ws.Range("G9").Value = 31.97572966389053
$ws.Range("H9").Value = 16.24649917737122
$ws.Range("I9").Value = 34.03823530509818
$ws.Range("J9").Value = 32.88807817580295
$ws.Range("K9").Value = 519.9884421999199
$ws.Range("L9").Value = 591.2424668176375
$ws.Range("M9").Value = 1189.880043115714
$ws.Range("N9").Value = 32693.4987829374
$ws.Range("O9").Value = 0.495
$ws.Range("P9").Value = 4.867736679452733
$ws.Range("Q9").Value = 10.00650944243258
$ws.Range("R9").Value = 0.803580671269581
$ws.Range("S9").Value = 0.426565594452296
$ws.Range("T9").Value = 1.0540319915
$ws.Range("U9").Value = 23.8170023194615
$ws.Range("V9").Value = 84918180.47217181
$ws.Range("W9").Value = 23.063378380558
$ws.Range("Y9").Value = 0.03783119916848355
$ws.Range("AC9").Value = 0.615206496
$ws.Range("AG9").Value = 3.546834268193857

# Row 10
$ws.Range("G10").Value = 47.096738630358
$ws.Range("H10").Value = 34.0107282995287
$ws.Range("I10").Value = 53.0981788198579
$ws.Range("J10").Value = 51.40291339207745
$ws.Range("K10").Value = 1845.638977786491
$ws.Range("L10").Value = 1467.331118768128
$ws.Range("M10").Value = 2341.826300347025
$ws.Range("N10").Value = 104465.4086735491
$ws.Range("O10").Value = 0.965
$ws.Range("P10").Value = 9.697142391269981
$ws.Range("Q10").Value = 15.106813151
$ws.Range("R10").Value = 1.09679414795852
$ws.Range("S10").Value = 0.676093631735514
$ws.Range("T10").Value = 3.758692792
$ws.Range("U10").Value = 32.7991050880754
$ws.Range("V10").Value = 179871657.957079
$ws.Range("W10").Value = 31.65455160268605
$ws.Range("Y10").Value = 0.06619530851054861
$ws.Range("AC10").Value = 1.2935617725
$ws.Range("AG10").Value = 8.17129942893974

# Row 11
$ws.Range("G11").Value = 60.7081791426266
$ws.Range("H11").Value = 56.03985806408188
$ws.Range("I11").Value = 67.58796719799679
$ws.Range("J11").Value = 68.65201237478766
$ws.Range("K11").Value = 3562.589637956335
$ws.Range("L11").Value = 3044.774392097399
$ws.Range("M11").Value = 4008.670927503015
$ws.Range("N11").Value = 223742.2252884099
$ws.Range("O11").Value = 2.260976597
$ws.Range("P11").Value = 19.90659784617232
$ws.Range("Q11").Value = 26.433341954
$ws.Range("R11").Value = 1.57686094230327
$ws.Range("S11").Value = 1.058069735228585
$ws.Range("T11").Value = 8.55328689645
$ws.Range("U11").Value = 44.1439925015829
$ws.Range("V11").Value = 704706269.1745907
$ws.Range("W11").Value = 43.3666275645425
$ws.Range("Y11").Value = 0.09993507462686579
$ws.Range("AC11").Value = 2.764510773803105
$ws.Range("AG11").Value = 19.19298651717195

# Row 12
$ws.Range("G12").Value = 90.8127810177044
$ws.Range("H12").Value = 91.1906848126461
$ws.Range("I12").Value = 97.48458318188359
$ws.Range("J12").Value = 98.7431693989071
$ws.Range("K12").Value = 8889.461961586105
$ws.Range("L12").Value = 7544.160584938694
$ws.Range("M12").Value = 8170.701187797568
$ws.Range("N12").Value = 748929.4811692492
$ws.Range("O12").Value = 147.4
$ws.Range("P12").Value = 269091.608902911
$ws.Range("Q12").Value = 1460.526315789
$ws.Range("R12").Value = 5.35324488052067
$ws.Range("S12").Value = 4.48865892578156
$ws.Range("T12").Value = 84.6368217716
$ws.Range("U12").Value = 227.052991658391
$ws.Range("V12").Value = 26130292711.2759
$ws.Range("W12").Value = 259.160835724776
$ws.Range("Y12").Value = 139.4757719298264
$ws.Range("AC12").Value = 58102.297777778
$ws.Range("AG12").Value = 384416.5841470158
